$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "species" column header in I1
$ws.Range("I1").Value = "species"

# Fill species value for every data row (2-60)
$species = "A. elegantissima x B. muscatinei"
for ($r = 2; $r -le 60; $r++) {
    $ws.Cells.Item($r, 9).Value = $species
}

# Scroll the view so that row 35 is the top-left visible row, matching the
# saved sheetView topLeftCell="A35"
$excel.ActiveWindow.ScrollRow = 35
